$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 437: Friday separator row (copy style/format from an existing separator row) ---
$ws.Range("A69:F69").Copy($ws.Range("A437:F437"))
$ws.Cells.Item(437, 3).Value = "FRIDAY"

# --- Row 438: Pickup Mic entry (long wrapped note -> row height 60) ---
$ws.Range("A75:F75").Copy($ws.Range("A438:F438"))
$ws.Cells.Item(438, 1).Value = "Pickup Mic"
$ws.Cells.Item(438, 2).Value = 42671
$ws.Cells.Item(438, 3).Value = "1730"
$ws.Cells.Item(438, 4).Value = "VH"
$ws.Cells.Item(438, 5).Value = "A"
$ws.Cells.Item(438, 6).Value = "1 lec mic, 2 desk mics, 3 wireless Shure mics, stands, cables, audio mixer, receivers etc. on a cart.  Make certain to power off Shure wireless mics. Return all items to Vari 1019 MCR. Raul needs some of this equipment for Saturday in his area."
$ws.Rows.Item(438).RowHeight = 60

# --- Row 439: Monday separator row ---
$ws.Range("A69:F69").Copy($ws.Range("A439:F439"))
$ws.Cells.Item(439, 3).Value = "MONDAY"

# --- Row 440: Pickup Mic entry (row height 30) ---
$ws.Range("A407:F407").Copy($ws.Range("A440:F440"))
$ws.Cells.Item(440, 1).Value = "Pickup Mic"
$ws.Cells.Item(440, 2).Value = 42674
$ws.Cells.Item(440, 3).Value = "1600"
$ws.Cells.Item(440, 4).Value = "FC"
$ws.Cells.Item(440, 5).Value = "152-A.Hall"
$ws.Cells.Item(440, 6).Value = "Leave mic cables in place."
$ws.Rows.Item(440).RowHeight = 30

# --- Row 441: AV Shutdown entry ---
$ws.Range("A75:F75").Copy($ws.Range("A441:F441"))
$ws.Cells.Item(441, 1).Value = "AV Shutdown"
$ws.Cells.Item(441, 2).Value = 42674
$ws.Cells.Item(441, 3).Value = "1600"
$ws.Cells.Item(441, 4).Value = "FC"
$ws.Cells.Item(441, 5).Value = "152-A.Hall"
$ws.Cells.Item(441, 6).Value = "Return wireless keyboard & projector remote to FDRS 156A."

# --- Row 442: AV Shutdown entry ---
$ws.Range("A75:F75").Copy($ws.Range("A442:F442"))
$ws.Cells.Item(442, 1).Value = "AV Shutdown"
$ws.Cells.Item(442, 2).Value = 42674
$ws.Cells.Item(442, 3).Value = "1630"
$ws.Cells.Item(442, 4).Value = "MC"
$ws.Cells.Item(442, 5).Value = "101A"
$ws.Cells.Item(442, 6).Value = "Pick up wireless keyboard and TV remote control. To FDRS 164."

# --- Row 443: Pickup PC entry (row height 30) ---
$ws.Range("A72:F72").Copy($ws.Range("A443:F443"))
$ws.Cells.Item(443, 1).Value = "Pickup PC"
$ws.Cells.Item(443, 2).Value = 42674
$ws.Cells.Item(443, 3).Value = "1700"
$ws.Cells.Item(443, 4).Value = "VC"
$ws.Cells.Item(443, 5).Value = "010-SCR"
$ws.Cells.Item(443, 6).Value = "Leave portable screen and network cable, but get carpets. To  Vanier 040 storeroom."
$ws.Rows.Item(443).RowHeight = 30

# --- Row 444: Other / Door code entry ---
$ws.Range("A74:F74").Copy($ws.Range("A444:F444"))
$ws.Cells.Item(444, 1).Value = "Other"
$ws.Cells.Item(444, 2).Value = 42674
$ws.Cells.Item(444, 3).Value = "1730"
$ws.Cells.Item(444, 4).Value = "MC"
$ws.Cells.Item(444, 5).Value = "157A"
$ws.Cells.Item(444, 6).Value = "Door code"

# --- Row 445: AV Shutdown entry (no note) ---
$ws.Range("A162:E162").Copy($ws.Range("A445:E445"))
$ws.Cells.Item(445, 1).Value = "AV Shutdown"
$ws.Cells.Item(445, 2).Value = 42674
$ws.Cells.Item(445, 3).Value = "1730"
$ws.Cells.Item(445, 4).Value = "R"
$ws.Cells.Item(445, 5).Value = "N203"

# --- Row 446: AV Shutdown entry ---
$ws.Range("A75:F75").Copy($ws.Range("A446:F446"))
$ws.Cells.Item(446, 1).Value = "AV Shutdown"
$ws.Cells.Item(446, 2).Value = 42674
$ws.Cells.Item(446, 3).Value = "1800"
$ws.Cells.Item(446, 4).Value = "FC"
$ws.Cells.Item(446, 5).Value = "305"
$ws.Cells.Item(446, 6).Value = "Turn off PC and projector. Leave all in and lock room."

# --- Row 447: AV Shutdown entry (no note) ---
$ws.Range("A162:E162").Copy($ws.Range("A447:E447"))
$ws.Cells.Item(447, 1).Value = "AV Shutdown"
$ws.Cells.Item(447, 2).Value = 42674
$ws.Cells.Item(447, 3).Value = "1830"
$ws.Cells.Item(447, 4).Value = "R"
$ws.Cells.Item(447, 5).Value = "S203"

# --- Row 448: AV Shutdown entry ---
$ws.Range("A75:F75").Copy($ws.Range("A448:F448"))
$ws.Cells.Item(448, 1).Value = "AV Shutdown"
$ws.Cells.Item(448, 2).Value = 42674
$ws.Cells.Item(448, 3).Value = "2130"
$ws.Cells.Item(448, 4).Value = "R"
$ws.Cells.Item(448, 5).Value = "N102"
$ws.Cells.Item(448, 6).Value = "Nat Taylor Cinema. Lock cinema all doors after shutdown."

# --- Row 449: Pickup PC entry ---
$ws.Range("A75:F75").Copy($ws.Range("A449:F449"))
$ws.Cells.Item(449, 1).Value = "Pickup PC"
$ws.Cells.Item(449, 2).Value = 42674
$ws.Cells.Item(449, 3).Value = "2200"
$ws.Cells.Item(449, 4).Value = "MC"
$ws.Cells.Item(449, 5).Value = "140-SCR"
$ws.Cells.Item(449, 6).Value = "Door code 7083*. Leave portable screen. Return to FDRS 156A."

# --- Update the view state to match: frozen pane scrolled down, new selection ---
$ws.Range("C464").Select()
